$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 - SMOKE_ST
$ws.Range("F7").Value = "__BLANK__"
$ws.Range("G7").Value = "paste"
$ws.Range("H7").NumberFormat = "@"
$ws.Range("H7").Value = "4"
$ws.Range("I7").Value = "non-smokers as inclusion criterion, but we cannot be sure, if they are real never smokers or former smokers"
$ws.Range("J7").Value = "partial"
$ws.Range("K7").Value = "proximate"

# Row 10 - POTATOES_TUB_01
$ws.Range("I10").Value = "does not contain other tubers than potatoes"

# Row 13 - VEGETABLES_02
$ws.Range("F13").Value = "Gem02_1_NCI"

# Row 18 - MUSHROOMS_0205
$ws.Range("F18").Value = "Gem02_2_NCI"

# Row 23 - LEGUMES_TOT_03
$ws.Range("F23").Value = "Gem02_Hu_NCI"

# Row 25 - FRUITS_TOT_04
$ws.Range("F25").Value = "Obst_sum_NCI"

# Row 27 - NUTS_SEEDS_0402
$ws.Range("F27").Value = "Nusa_sum_NCI"

# Row 30 - DAIRY_05
$ws.Range("F30").Value = "Milc_sum_NCI"

# Row 31 - MILK_0501
$ws.Range("F31").Value = "milk_NCI"
$ws.Range("G31").Value = "direct_mapping"
$ws.Range("H31").Value = "direct_mapping"
$ws.Range("I31").Value = ""

# Row 32 - MILKBEV_0502
$ws.Range("F32").Value = "milkbased_bev_NCI"
$ws.Range("G32").Value = "direct_mapping"
$ws.Range("H32").Value = "direct_mapping"
$ws.Range("I32").Value = ""
$ws.Range("J32").Value = "complete"
$ws.Range("K32").Value = "identical"

# Row 34 - CURD_0504
$ws.Range("F34").Value = "quark_curd_NCI"
$ws.Range("G34").Value = "direct_mapping"
$ws.Range("H34").Value = "direct_mapping"
$ws.Range("I34").Value = ""

# Row 35 - CHEESE_0505
$ws.Range("F35").Value = "cheeses_NCI"
$ws.Range("G35").Value = "direct_mapping"
$ws.Range("H35").Value = "direct_mapping"
$ws.Range("I35").Value = ""

# Row 36 - DAIRYDESSERT_0506
$ws.Range("F36").Value = "Sues13_6_NCI"

# Row 41 - CEREAL_PROD_06
$ws.Range("F41").Value = "Brot_sum_NCI"

# Row 44 - BREAD_PROD_0603
$ws.Range("F44").Value = "Brot01_1_NCI"

# Row 50 - MEAT_PROD_07
$ws.Range("F50").Value = "Flei_sum_NCI"

# Row 66 - PROCMEAT_0704
$ws.Range("F66").Value = "Flei_Wurst_NCI"

# Row 68 - FISH_SHELLFISH_08
$ws.Range("F68").Value = "Fish_sum_NCI"

# Row 74 - FAT_10
$ws.Range("F74").Value = "Fett_sum_NCI"

# Row 81 - SUGAR_CONFECT_11
$ws.Range("F81").Value = "Sues_sum_NCI"

# Row 83 - CHOCOLATE_1102
$ws.Range("F83").Value = "chocolate_sweets_NCI"
$ws.Range("G83").Value = "direct_mapping"
$ws.Range("H83").Value = "direct_mapping"
$ws.Range("J83").Value = "complete"
$ws.Range("K83").Value = "identical"

# Row 84 - NONCHOC_SWEETS_1103
$ws.Range("F84").Value = "nonchoc_sweets_NCI"
$ws.Range("G84").Value = "direct_mapping"
$ws.Range("H84").Value = "direct_mapping"
$ws.Range("J84").Value = "complete"
$ws.Range("K84").Value = "identical"

# Row 86 - ICECREAM_1105
$ws.Range("F86").Value = "Sues13_2_NCI"

# Row 93 - NONALC_BEV_13
$ws.Range("F93").Value = "Getr15_1_NCI"

# Row 94 - FRUITVEG_JUICE_1301
$ws.Range("F94").Value = "Getr15_14_NCI"

# Row 95 - SOFTDRINKS_1302
$ws.Range("F95").Value = "Getr15_16_NCI"

# Row 96 - HOTDRINKS_1303
$ws.Range("F96").Value = "Kaffee_NCI;blackgreentea_NCI;Getr15_13_NCI"
$ws.Range("H96").Value = "Kaffee_NCI+blackgreentea_NCI+Getr15_13_NCI"

# Row 97 - COFFEE_130301
$ws.Range("F97").Value = "Kaffee_NCI"
$ws.Range("G97").Value = "direct_mapping"
$ws.Range("H97").Value = "direct_mapping"
$ws.Range("I97").Value = ""
$ws.Range("J97").Value = "complete"
$ws.Range("K97").Value = "identical"

# Row 98 - TEA_130302
$ws.Range("F98").Value = "blackgreentea_NCI"
$ws.Range("G98").Value = "direct_mapping"
$ws.Range("H98").Value = "direct_mapping"
$ws.Range("I98").Value = ""
$ws.Range("J98").Value = "complete"
$ws.Range("K98").Value = "identical"

# Row 99 - HERBALTEA_130303
$ws.Range("F99").Value = "Getr15_13_NCI"

# Row 101 - WATER_1304
$ws.Range("F101").Value = "Getr15_11_NCI"

# Row 102 - ALC_BEV_14
$ws.Range("F102").Value = "Getr15_2_NCI"

# Row 103 - WINE_1401
$ws.Range("F103").Value = "Getr15_22_NCI"

# Row 105 - BEER_1403
$ws.Range("F105").Value = "Getr15_21_NCI"

# Row 106 - SPIRITS_1404
$ws.Range("F106").Value = "Getr15_23_MW"
$ws.Range("I106").Value = "calculated as mean of 2 24h-recalls"
$ws.Range("J106").Value = "partial"
$ws.Range("K106").Value = "tentative"

# Row 109 - COCKTAILS_1407
$ws.Range("F109").Value = "Getr15_24_MW"
$ws.Range("I109").Value = "calculated as mean of 2 24h-recalls"
$ws.Range("J109").Value = "partial"
$ws.Range("K109").Value = "tentative"

# Row 110 - CONDIMENT_SAUCES_15
$ws.Range("F110").Value = "Sose_sum_NCI"

# Row 119 - SOUP_BOUILLON_16
$ws.Range("F119").Value = "Supp_sum_NCI"

# Row 122 - MISCELLANEOUS_17
$ws.Range("F122").Value = "Sonst_sum_NCI"

# Row 126 - ART_SWEETENER_170201
$ws.Range("F126").Value = "Suessstoffe_MW"
$ws.Range("I126").Value = "calculated as mean of 2 24h-recalls"
$ws.Range("J126").Value = "partial"
$ws.Range("K126").Value = "tentative"

# Row 129 - DIETARY_ASSESS_INSTR
$ws.Range("I129").Value = "2 (24HDR) NCI method"
